# Auto-generated edit script: updates market-price-derived columns (H-N)
# across the 8 crafting-class sheets, per scheduled runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 871.4286
$ws.Range("I58").Value = 100
$ws.Range("K58").Value = 300
$ws.Range("M58").Value = -150
$ws.Range("H94").Value = 414.125
$ws.Range("I94").Value = 414.125
$ws.Range("K94").Value = 414.125
$ws.Range("M94").Value = 36.875
$ws.Range("H112").Value = 1919.0968
$ws.Range("I112").Value = 1000
$ws.Range("J112").Value = 1949.7333
$ws.Range("K112").Value = 3000
$ws.Range("L112").Value = 5849.199900000001
$ws.Range("M112").Value = -1892
$ws.Range("N112").Value = -8065.199900000001
$ws.Range("H132").Value = 18544.861
$ws.Range("I132").Value = 1361.24
$ws.Range("J132").Value = 125942.5
$ws.Range("K132").Value = 4083.72
$ws.Range("L132").Value = 377827.5
$ws.Range("M132").Value = -1553.72
$ws.Range("N132").Value = -382887.5
$ws.Range("H138").Value = 2337.3
$ws.Range("I138").Value = 1316.8572
$ws.Range("K138").Value = 3950.5716
$ws.Range("M138").Value = 1189.4284

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3080.1428
$ws.Range("I61").Value = 1101.4286
$ws.Range("K61").Value = 1101.4286
$ws.Range("M61").Value = -889.4286
$ws.Range("H136").Value = 3080.1428
$ws.Range("I136").Value = 1101.4286
$ws.Range("K136").Value = 3304.2858
$ws.Range("M136").Value = -754.2857999999997

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3711.35
$ws.Range("I86").Value = 1651.5555
$ws.Range("K86").Value = 1651.5555
$ws.Range("M86").Value = -528.5554999999999
$ws.Range("H89").Value = 3711.35
$ws.Range("I89").Value = 1651.5555
$ws.Range("K89").Value = 8257.7775
$ws.Range("M89").Value = -2641.7775
$ws.Range("H94").Value = 5628.909
$ws.Range("I94").Value = 4988.5
$ws.Range("K94").Value = 4988.5
$ws.Range("M94").Value = -4537.5
$ws.Range("H134").Value = 3093.3242
$ws.Range("I134").Value = 1584.5652
$ws.Range("K134").Value = 4753.6956
$ws.Range("M134").Value = -2218.6956

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 434.44446
$ws.Range("I22").Value = 451.125
$ws.Range("J22").Value = 301
$ws.Range("K22").Value = 451.125
$ws.Range("L22").Value = 301
$ws.Range("M22").Value = -101.125
$ws.Range("N22").Value = -1001
$ws.Range("H31").Value = 3873.2856
$ws.Range("J31").Value = 4724.6
$ws.Range("L31").Value = 4724.6
$ws.Range("N31").Value = -5314.6
$ws.Range("H32").Value = 4600
$ws.Range("J32").Value = 4600
$ws.Range("L32").Value = 4600
$ws.Range("N32").Value = -5232
$ws.Range("H34").Value = 3873.2856
$ws.Range("J34").Value = 4724.6
$ws.Range("L34").Value = 4724.6
$ws.Range("N34").Value = -5128.6
$ws.Range("H43").Value = 18821.166
$ws.Range("J43").Value = 18821.166
$ws.Range("L43").Value = 18821.166
$ws.Range("N43").Value = -19189.166
$ws.Range("H74").Value = 49999
$ws.Range("I74").Value = 32000
$ws.Range("J74").Value = 55998.668
$ws.Range("K74").Value = 32000
$ws.Range("L74").Value = 55998.668
$ws.Range("M74").Value = -31126
$ws.Range("N74").Value = -57746.668
$ws.Range("H77").Value = 49999
$ws.Range("I77").Value = 32000
$ws.Range("J77").Value = 55998.668
$ws.Range("K77").Value = 96000
$ws.Range("L77").Value = 167996.004
$ws.Range("M77").Value = -91632
$ws.Range("N77").Value = -176732.004
$ws.Range("H101").Value = 18821.166
$ws.Range("J101").Value = 18821.166
$ws.Range("L101").Value = 18821.166
$ws.Range("N101").Value = -25311.166

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 64987.25
$ws.Range("J59").Value = 86316.336
$ws.Range("L59").Value = 258949.008
$ws.Range("N59").Value = -260029.008
$ws.Range("H63").Value = 172502.33
$ws.Range("I63").Value = 501888.75
$ws.Range("J63").Value = 7809.125
$ws.Range("K63").Value = 1505666.25
$ws.Range("L63").Value = 23427.375
$ws.Range("M63").Value = -1504917.25
$ws.Range("N63").Value = -24925.375
$ws.Range("H66").Value = 172502.33
$ws.Range("I66").Value = 501888.75
$ws.Range("J66").Value = 7809.125
$ws.Range("K66").Value = 4516998.75
$ws.Range("L66").Value = 70282.125
$ws.Range("M66").Value = -4513254.75
$ws.Range("N66").Value = -77770.125
$ws.Range("H103").Value = 4534.857
$ws.Range("I103").Value = 4497.25
$ws.Range("J103").Value = 4585
$ws.Range("K103").Value = 13491.75
$ws.Range("L103").Value = 13755
$ws.Range("M103").Value = -12612.75
$ws.Range("N103").Value = -15513
$ws.Range("H121").Value = 2081.0344
$ws.Range("I121").Value = 235.83333
$ws.Range("J121").Value = 2562.3914
$ws.Range("K121").Value = 707.49999
$ws.Range("L121").Value = 7687.174199999999
$ws.Range("M121").Value = 602.50001
$ws.Range("N121").Value = -10307.1742

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4356.8076
$ws.Range("I132").Value = 4193.4287
$ws.Range("J132").Value = 5043
$ws.Range("K132").Value = 12580.2861
$ws.Range("L132").Value = 15129
$ws.Range("M132").Value = -10050.2861
$ws.Range("N132").Value = -20189

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4889.6924
$ws.Range("I7").Value = 2692.2
$ws.Range("K7").Value = 2692.2
$ws.Range("M7").Value = -2580.2
$ws.Range("H40").Value = 10590.682
$ws.Range("I40").Value = 13981.777
$ws.Range("K40").Value = 13981.777
$ws.Range("M40").Value = -13845.777
$ws.Range("H100").Value = 95591.586
$ws.Range("I100").Value = 112570
$ws.Range("J100").Value = 10699.5
$ws.Range("K100").Value = 112570
$ws.Range("L100").Value = 10699.5
$ws.Range("M100").Value = -112029
$ws.Range("N100").Value = -11781.5
$ws.Range("H102").Value = 34199.2
$ws.Range("J102").Value = 34199.2
$ws.Range("L102").Value = 34199.2
$ws.Range("N102").Value = -40689.2
$ws.Range("H105").Value = 39307.5
$ws.Range("J105").Value = 39307.5
$ws.Range("L105").Value = 39307.5
$ws.Range("N105").Value = -46295.5
$ws.Range("H126").Value = 4889.6924
$ws.Range("I126").Value = 2692.2
$ws.Range("K126").Value = 8076.599999999999
$ws.Range("M126").Value = -5606.599999999999
$ws.Range("H132").Value = 4861.64
$ws.Range("I132").Value = 4132.5557
$ws.Range("J132").Value = 6736.4287
$ws.Range("K132").Value = 12397.6671
$ws.Range("L132").Value = 20209.2861
$ws.Range("M132").Value = -9867.667099999999
$ws.Range("N132").Value = -25269.2861
$ws.Range("H136").Value = 5889.852
$ws.Range("I136").Value = 4101.7334
$ws.Range("J136").Value = 8125
$ws.Range("K136").Value = 12305.2002
$ws.Range("L136").Value = 24375
$ws.Range("M136").Value = -9755.200199999999
$ws.Range("N136").Value = -29475

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 12860.25
$ws.Range("I52").Value = 3648.6667
$ws.Range("J52").Value = 40495
$ws.Range("K52").Value = 3648.6667
$ws.Range("L52").Value = 40495
$ws.Range("M52").Value = -3422.6667
$ws.Range("N52").Value = -40947
$ws.Range("H64").Value = 59957
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 59957
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 59957
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -60453
$ws.Range("H67").Value = 59957
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 59957
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 59957
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -61673
$ws.Range("H132").Value = 4062.8667
$ws.Range("I132").Value = 3746.56
$ws.Range("K132").Value = 11239.68
$ws.Range("M132").Value = -8709.68
$ws.Range("H136").Value = 5902.15
$ws.Range("I136").Value = 5092.357
$ws.Range("K136").Value = 15277.071
$ws.Range("M136").Value = -12727.071
$ws.Range("H139").Value = 70000
$ws.Range("J139").Value = 70000
$ws.Range("L139").Value = 70000
$ws.Range("N139").Value = -80280
